$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# NOTE: '83÷3=' is the target of one replacement ('21÷9=' -> '83÷3=') and the
# source of another ('83÷3=' -> '56÷6='). To avoid the newly-inserted '83÷3='
# being matched by the later replacement, the '83÷3=' -> '56÷6=' step is
# performed before '21÷9=' -> '83÷3=' is applied.

Replace-Text "24÷5=" "50÷3="
Replace-Text "41÷2=" "99÷3="
Replace-Text "75÷9=" "89÷4="
Replace-Text "68÷8=" "42÷2="
Replace-Text "32÷8=" "65÷8="

Replace-Text "83÷3=" "56÷6="
Replace-Text "21÷9=" "83÷3="

Replace-Text "33÷9=" "54÷6="
Replace-Text "77÷5=" "60÷5="
Replace-Text "45÷8=" "97÷9="
Replace-Text "72÷5=" "50÷2="
Replace-Text "67÷9=" "81÷8="
Replace-Text "68÷7=" "32÷2="
Replace-Text "49÷7=" "61÷5="
Replace-Text "50÷5=" "19÷8="
Replace-Text "10÷6=" "58÷7="
Replace-Text "62÷4=" "33÷4="
Replace-Text "20÷3=" "19÷9="
Replace-Text "67÷7=" "46÷3="
Replace-Text "97÷3=" "25÷6="
Replace-Text "62÷8=" "67÷4="
Replace-Text "69÷8=" "27÷2="
Replace-Text "56÷3=" "46÷2="
Replace-Text "85÷7=" "47÷5="
Replace-Text "82÷9=" "11÷7="
